# "working on analysis of Exp3"
# - D13: N Total Surveys for the Adapt_Emo_Identity_JAPMEPS (baseline-removed) row bumps 47 -> 48
# - D15 (Missings row): was a plain "1", now annotated as text "3 (wg missings)"
# - D16 (ToDo row): was red "5 -> check", now regular-colored "3 removed, 2 kept in "
#   (matches the formatting already used by the neighboring C16 cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric update
$ws.Range("D13").Value = 48

# D16 gets new text content first (while it still carries its own red style),
# then D15 is written - this preserves the shared-string insertion order seen
# in the saved file (index 63 = "3 removed, 2 kept in ", 64 = "3 (wg missings)").
$ws.Range("D16").Value = "3 removed, 2 kept in "
$ws.Range("D15").Value = "3 (wg missings)"

# Re-format D16 to match C16 (drop the red "check me" font, use the regular one)
$ws.Range("C16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null

# Move the active selection to D32, as recorded in the saved view state
$ws.Range("D32").Select() | Out-Null
